# Update cryptocurrency price/volume figures per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '27.407.51'
$ws.Range('E2').Value = '  -0.95%  '
$ws.Range('D3').Value = '1.566.42'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '''208.68'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('D6').Value = '''0.501'
$ws.Range('E6').Value = '  -0.74%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '''22.04'
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '1.789.80'
$ws.Range('D13').Value = '1.569.26'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('E15').Value = '  -2.82%  '
$ws.Range('D16').Value = '''63.53'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').Value = '27.424.25'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '''212.90'
$ws.Range('E18').Value = '  -3.35%  '
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '''7.26'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '''4.11'
$ws.Range('E22').Value = '  -0.77%  '
$ws.Range('D23').Value = '''9.53'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  +2.49%  '
$ws.Range('D25').Value = '''152.88'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').Value = '''6.70'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').Value = '''14.97'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').Value = '''1.15'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = '1.374.54'
$ws.Range('E33').Value = '  -0.83%  '
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('D36').Value = '''0.965'
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').Value = '''0.0167'
$ws.Range('E38').Value = '  +1.19%  '
$ws.Range('D39').Value = '''0.531'
$ws.Range('E39').Value = '  -2.00%  '
$ws.Range('D40').Value = '''0.822'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').Value = '''0.974'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('D43').Value = '''1.79'
$ws.Range('E43').Value = '  +0.45%  '
$ws.Range('D44').Value = '''63.93'
$ws.Range('E44').Value = '  +0.29%  '
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').Value = '1.702.08'
$ws.Range('E47').Value = '  -1.16%  '
$ws.Range('D48').Value = '''85.47'
$ws.Range('E48').Value = '  -2.88%  '
$ws.Range('E49').Value = '  -1.46%  '
$ws.Range('E50').Value = '  -1.53%  '
$ws.Range('E51').Value = '  -0.95%  '
